# Update transition-probability matrix values on Sheet1 (North Ala._B)
# Values recomputed after simulating more games (updated Markov transition matrix).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2205128205128205
$ws.Range("C2").Value = 0.4948717948717949
$ws.Range("J2").Value = 0.01282051282051282
$ws.Range("P2").Value = 0.1897435897435897
$ws.Range("S2").Value = 0.08205128205128205

# Row 3
$ws.Range("B3").Value = 0.02475247524752475
$ws.Range("C3").Value = 0.0396039603960396
$ws.Range("J3").Value = 0.0396039603960396
$ws.Range("P3").Value = 0.7128712871287128
$ws.Range("S3").Value = 0.1831683168316832

# Row 4
$ws.Range("J4").Value = 0.01694915254237288
$ws.Range("P4").Value = 0.7796610169491526
$ws.Range("S4").Value = 0.2033898305084746

# Row 5
$ws.Range("P5").Value = 1

# Row 6
$ws.Range("B6").Value = 0.095
$ws.Range("D6").Value = 0.015
$ws.Range("F6").Value = 0.07000000000000001
$ws.Range("J6").Value = 0.22
$ws.Range("O6").Value = 0.015
$ws.Range("Q6").Value = 0.16
$ws.Range("R6").Value = 0.075
$ws.Range("S6").Value = 0.35

# Row 7
$ws.Range("B7").Value = 0.1481481481481481
$ws.Range("D7").Value = 0.03240740740740741
$ws.Range("F7").Value = 0.04629629629629629
$ws.Range("J7").Value = 0.1435185185185185
$ws.Range("O7").Value = 0.01388888888888889
$ws.Range("Q7").Value = 0.125
$ws.Range("R7").Value = 0.09259259259259259
$ws.Range("S7").Value = 0.3981481481481481

# Row 8
$ws.Range("B8").Value = 0.1257606490872211
$ws.Range("D8").Value = 0.03042596348884381
$ws.Range("E8").Value = 0.004056795131845842
$ws.Range("F8").Value = 0.05882352941176471
$ws.Range("J8").Value = 0.1095334685598377
$ws.Range("O8").Value = 0.02434077079107505
$ws.Range("Q8").Value = 0.1744421906693712
$ws.Range("R8").Value = 0.09127789046653144
$ws.Range("S8").Value = 0.3813387423935091

# Row 9
$ws.Range("B9").Value = 0.1414634146341463
$ws.Range("D9").Value = 0.02926829268292683
$ws.Range("F9").Value = 0.08780487804878048
$ws.Range("J9").Value = 0.07317073170731707
$ws.Range("O9").Value = 0.02439024390243903
$ws.Range("Q9").Value = 0.1560975609756098
$ws.Range("R9").Value = 0.08292682926829269
$ws.Range("S9").Value = 0.4048780487804878

# Row 10
$ws.Range("B10").Value = 0.1246031746031746
$ws.Range("D10").Value = 0.02222222222222222
$ws.Range("E10").Value = 0.0007936507936507937
$ws.Range("F10").Value = 0.05952380952380952
$ws.Range("J10").Value = 0.119047619047619
$ws.Range("O10").Value = 0.0253968253968254
$ws.Range("Q10").Value = 0.2261904761904762
$ws.Range("R10").Value = 0.07698412698412699
$ws.Range("S10").Value = 0.3452380952380952

# Row 11
$ws.Range("G11").Value = 0.1508196721311476
$ws.Range("J11").Value = 0.07868852459016394
$ws.Range("K11").Value = 0.1836065573770492
$ws.Range("L11").Value = 0.5836065573770491
$ws.Range("S11").Value = 0.003278688524590164

# Row 12
$ws.Range("G12").Value = 0.7914438502673797
$ws.Range("J12").Value = 0.1443850267379679
$ws.Range("K12").Value = 0.0053475935828877
$ws.Range("L12").Value = 0.03208556149732621
$ws.Range("S12").Value = 0.0267379679144385

# Row 13
$ws.Range("G13").Value = 0.6888888888888889
$ws.Range("J13").Value = 0.2888888888888889
$ws.Range("S13").Value = 0.02222222222222222

# Row 14
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5

# Row 15
$ws.Range("F15").Value = 0.008888888888888889
$ws.Range("H15").Value = 0.1555555555555556
$ws.Range("I15").Value = 0.08444444444444445
$ws.Range("J15").Value = 0.3244444444444444
$ws.Range("K15").Value = 0.05777777777777778
$ws.Range("M15").Value = 0.01777777777777778
$ws.Range("O15").Value = 0.05333333333333334
$ws.Range("S15").Value = 0.2977777777777778

# Row 16
$ws.Range("F16").Value = 0.003891050583657588
$ws.Range("H16").Value = 0.1945525291828794
$ws.Range("I16").Value = 0.09727626459143969
$ws.Range("J16").Value = 0.3618677042801556
$ws.Range("K16").Value = 0.1245136186770428
$ws.Range("M16").Value = 0.01167315175097276
$ws.Range("O16").Value = 0.0622568093385214
$ws.Range("S16").Value = 0.1439688715953307

# Row 17
$ws.Range("F17").Value = 0.01731601731601732
$ws.Range("H17").Value = 0.20995670995671
$ws.Range("I17").Value = 0.0735930735930736
$ws.Range("J17").Value = 0.4199134199134199
$ws.Range("K17").Value = 0.09740259740259741
$ws.Range("M17").Value = 0.01515151515151515
$ws.Range("N17").Value = 0.004329004329004329
$ws.Range("O17").Value = 0.05194805194805195
$ws.Range("S17").Value = 0.1103896103896104

# Row 18
$ws.Range("F18").Value = 0.0155440414507772
$ws.Range("H18").Value = 0.2435233160621762
$ws.Range("I18").Value = 0.09326424870466321
$ws.Range("J18").Value = 0.4404145077720207
$ws.Range("K18").Value = 0.08808290155440414
$ws.Range("M18").Value = 0.01036269430051814
$ws.Range("O18").Value = 0.0155440414507772
$ws.Range("S18").Value = 0.09326424870466321

# Row 19
$ws.Range("F19").Value = 0.01354581673306773
$ws.Range("H19").Value = 0.2127490039840637
$ws.Range("I19").Value = 0.08685258964143426
$ws.Range("J19").Value = 0.3657370517928287
$ws.Range("K19").Value = 0.1083665338645418
$ws.Range("M19").Value = 0.02310756972111554
$ws.Range("N19").Value = 0.001593625498007968
$ws.Range("O19").Value = 0.06772908366533864
$ws.Range("S19").Value = 0.1203187250996016
